$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-11 with new cluster pairing / TPM values
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Ret"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.534048
$ws.Range("H2").Value = 1.602144
$ws.Range("I2").Value = 0.2492808729834395
$ws.Range("J2").Value = 0.3324807621550537
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8639135
$ws.Range("N2").Value = 1.727827
$ws.Range("O2").Value = 0.1895490737713731
$ws.Range("P2").Value = 0.1380033232738433
$ws.Range("Q2").Value = 0.4613712768479999
$ws.Range("R2").Value = 2.768227661088
$ws.Range("S2").Value = 0.04725095858293024
$ws.Range("T2").Value = 0.04588345010201768

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Ret"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.534048
$ws.Range("H3").Value = 1.602144
$ws.Range("I3").Value = 0.2492808729834395
$ws.Range("J3").Value = 0.3324807621550537
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.384385
$ws.Range("N3").Value = 10.153155
$ws.Range("O3").Value = 0.7425593442349591
$ws.Range("P3").Value = 0.8109429541930055
$ws.Range("Q3").Value = 1.80742404048
$ws.Range("R3").Value = 16.26681636432
$ws.Range("S3").Value = 0.1851058415729009
$ws.Range("T3").Value = 0.2696229314743613

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Ret"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.534048
$ws.Range("H4").Value = 1.602144
$ws.Range("I4").Value = 0.2492808729834395
$ws.Range("J4").Value = 0.3324807621550537
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2890925
$ws.Range("N4").Value = 0.578185
$ws.Range("O4").Value = 0.06342905349812297
$ws.Range("P4").Value = 0.04618023185601746
$ws.Range("Q4").Value = 0.15438927144
$ws.Range("R4").Value = 0.92633562864
$ws.Range("S4").Value = 0.01581164982852538
$ws.Range("T4").Value = 0.01535403868398577

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Ret"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.534048
$ws.Range("H5").Value = 1.602144
$ws.Range("I5").Value = 0.2492808729834395
$ws.Range("J5").Value = 0.3324807621550537
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01599833333333333
$ws.Range("N5").Value = 0.047995
$ws.Range("O5").Value = 0.003510153811948785
$ws.Range("P5").Value = 0.003833410116017465
$ws.Range("Q5").Value = 0.00854387792
$ws.Range("R5").Value = 0.07689490128000001
$ws.Range("S5").Value = 0.0008750142065487409
$ws.Range("T5").Value = 0.00127453511702638

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Ret"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.534048
$ws.Range("H6").Value = 1.602144
$ws.Range("I6").Value = 0.2492808729834395
$ws.Range("J6").Value = 0.3324807621550537
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.004340666666666667
$ws.Range("N6").Value = 0.013022
$ws.Range("O6").Value = 0.0009523746835961471
$ws.Range("P6").Value = 0.001040080561116354
$ws.Range("Q6").Value = 0.002318124352
$ws.Range("R6").Value = 0.020863119168
$ws.Range("S6").Value = 0.0002374087925341745
$ws.Range("T6").Value = 0.0003458067776626214

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Ret"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.6083065
$ws.Range("H7").Value = 3.216613
$ws.Range("I7").Value = 0.7507191270165605
$ws.Range("J7").Value = 0.6675192378449464
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8639135
$ws.Range("N7").Value = 1.727827
$ws.Range("O7").Value = 0.1895490737713731
$ws.Range("P7").Value = 0.1380033232738433
$ws.Range("Q7").Value = 1.38943769748775
$ws.Range("R7").Value = 5.557750789951
$ws.Range("S7").Value = 0.1422981151884428
$ws.Range("T7").Value = 0.09211987317182563

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Artn"
$ws.Range("C8").Value = "Ret"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.6083065
$ws.Range("H8").Value = 3.216613
$ws.Range("I8").Value = 0.7507191270165605
$ws.Range("J8").Value = 0.6675192378449464
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.384385
$ws.Range("N8").Value = 10.153155
$ws.Range("O8").Value = 0.7425593442349591
$ws.Range("P8").Value = 0.8109429541930055
$ws.Range("Q8").Value = 5.4431283940025
$ws.Range("R8").Value = 32.658770364015
$ws.Range("S8").Value = 0.5574535026620581
$ws.Range("T8").Value = 0.5413200227186443

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Artn"
$ws.Range("C9").Value = "Ret"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.6083065
$ws.Range("H9").Value = 3.216613
$ws.Range("I9").Value = 0.7507191270165605
$ws.Range("J9").Value = 0.6675192378449464
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2890925
$ws.Range("N9").Value = 0.578185
$ws.Range("O9").Value = 0.06342905349812297
$ws.Range("P9").Value = 0.04618023185601746
$ws.Range("Q9").Value = 0.46494934685125
$ws.Range("R9").Value = 1.859797387405
$ws.Range("S9").Value = 0.0476174036695976
$ws.Range("T9").Value = 0.03082619317203169

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Artn"
$ws.Range("C10").Value = "Ret"
$ws.Range("D10").Value = "Neutrophils"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.6083065
$ws.Range("H10").Value = 3.216613
$ws.Range("I10").Value = 0.7507191270165605
$ws.Range("J10").Value = 0.6675192378449464
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.01599833333333333
$ws.Range("N10").Value = 0.047995
$ws.Range("O10").Value = 0.003510153811948785
$ws.Range("P10").Value = 0.003833410116017465
$ws.Range("Q10").Value = 0.02573022348916667
$ws.Range("R10").Value = 0.154381340935
$ws.Range("S10").Value = 0.002635139605400044
$ws.Range("T10").Value = 0.002558874998991086

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Artn"
$ws.Range("C11").Value = "Ret"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.6083065
$ws.Range("H11").Value = 3.216613
$ws.Range("I11").Value = 0.7507191270165605
$ws.Range("J11").Value = 0.6675192378449464
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.004340666666666667
$ws.Range("N11").Value = 0.013022
$ws.Range("O11").Value = 0.0009523746835961471
$ws.Range("P11").Value = 0.001040080561116354
$ws.Range("Q11").Value = 0.006981122414333334
$ws.Range("R11").Value = 0.04188673448600001
$ws.Range("S11").Value = 0.0007149658910619726
$ws.Range("T11").Value = 0.0006942737834537331

# Remove old rows 12 and 13 which no longer exist in the new data (13 rows -> 11 rows)
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(12).Delete()
